# This workbook contains a weekly fruit/vegetable price log ("Uva" - grapes)
# sheet where each row is one price observation. A new weekly observation is
# inserted at row 93 (pushing all the existing rows from 93 downward by one),
# and it is populated with a new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 93; this shifts rows 93:118
# down to 94:119, matching the target layout (dimension becomes A1:T119).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly price observation.
# Single-quoted strings are used throughout so values such as "$/caja 20
# kilos" or accented names are taken literally (no PowerShell interpolation).
$ws.Cells.Item(93, 1).Value  = 1
$ws.Cells.Item(93, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(93, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(93, 4).Value  = (Get-Date -Year 2022 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(93, 5).Value  = 15
$ws.Cells.Item(93, 6).Value  = 'Fruta'
$ws.Cells.Item(93, 7).Value  = 100109
$ws.Cells.Item(93, 8).Value  = 'Uva'
$ws.Cells.Item(93, 9).Value  = 100109001
$ws.Cells.Item(93, 10).Value = 'Uva'
$ws.Cells.Item(93, 11).Value = 'Superior Seedless'
$ws.Cells.Item(93, 12).Value = 'Primera'
$ws.Cells.Item(93, 13).Value = 450
$ws.Cells.Item(93, 14).Value = 27000
$ws.Cells.Item(93, 15).Value = 28000
$ws.Cells.Item(93, 16).Value = 27556
$ws.Cells.Item(93, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(93, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(93, 19).Value = 1378
$ws.Cells.Item(93, 20).Value = 20
